$d = $word.ActiveDocument

# --- Change 1: the paragraph holding "{ some code, part of Definition 2 }"
# has two conflicting <w:pStyle> entries (SourceCode + Definition). Re-apply
# a single, unambiguous style ("Source Code") so the duplicate collapses.
$codePara = $null
For ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    If ($para.Range.Text -like "*some code, part of Definition 2*") {
        $codePara = $para
    }
}
If ($codePara -ne $null) {
    $codePara.Style = "Source Code"
}

# --- Change 2: "Block Text" should look like a normal block quote -
# indented left/right, using the same font/size as body text (i.e. no
# direct character formatting override), instead of the old unindented,
# smaller, different-typeface look.
$old = $d.Styles("BlockText")
$old.Delete()

$bt = $d.Styles.Add("Block Text", 1)
$bt.BaseStyle = $d.Styles("BodyText")
$bt.NextParagraphStyle = $d.Styles("BodyText")
$bt.Priority = 9
$bt.UnhideWhenUsed = $true
$bt.QuickStyle = $true

$bt.ParagraphFormat.SpaceBefore = 5
$bt.ParagraphFormat.SpaceAfter = 5
$bt.ParagraphFormat.FirstLineIndent = 0
$bt.ParagraphFormat.LeftIndent = 24
$bt.ParagraphFormat.RightIndent = 24
